$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 335 (shifts old rows 335-380 down to 336-381,
# carrying their existing values/formatting with them).
$ws.Rows.Item(335).Insert()

# Populate the newly inserted row 335 with the new weekly price entry.
$ws.Range("A335").Value = 9
$ws.Range("B335").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C335").Value = "Metropolitana"
$ws.Range("D335").Value = 45131
$ws.Range("D335").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E335").Value = 13
$ws.Range("F335").Value = 100112001
$ws.Range("G335").Value = "Berenjena"
$ws.Range("H335").Value = "Sin especificar"
$ws.Range("I335").Value = "Primera"
$ws.Range("J335").Value = 70
$ws.Range("K335").Value = 6000
$ws.Range("L335").Value = 7000
$ws.Range("M335").Value = 6500
$ws.Range("N335").Value = "$/caja 60 unidades"
$ws.Range("O335").Value = "Región de Arica y Parinacota"
$ws.Range("P335").Value = 108
$ws.Range("Q335").Value = 60
$ws.Range("R335").Value = "Hortaliza"
